# Complete SmartBlind testing on #21
# Fill in the test-completion Date (column C) and File Name (column E)
# for the SmartCurtain (CX), Routine (RX) and CloudSync (CS) sections
# of the Summary sheet, and move the active selection to E70.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$dateMar28 = "28-03-2025"
$dateMar26 = "26-03-2025"
$cloudSyncFile = "FW_Android_App_Test_Suite_v1.7.4_CloudSyn_28Mar2024"

# SmartCurtain (CX) test rows
$ws.Range("C43").Value = $dateMar28
$ws.Range("C44").Value = $dateMar28
$ws.Range("C45").Value = $dateMar28
$ws.Range("C46").Value = $dateMar28
$ws.Range("C48").Value = $dateMar28
$ws.Range("C49").Value = $dateMar28
$ws.Range("C51").Value = $dateMar28
$ws.Range("C52").Value = $dateMar28

# Routine (RX) test rows
$ws.Range("C66").Value = $dateMar26
$ws.Range("C68").Value = $dateMar26
$ws.Range("C69").Value = $dateMar26
$ws.Range("C70").Value = $dateMar26
$ws.Range("C71").Value = $dateMar26
$ws.Range("C72").Value = $dateMar26
$ws.Range("C73").Value = $dateMar26
$ws.Range("C74").Value = $dateMar26

# CloudSync (CS) test rows
$ws.Range("C77").Value = $dateMar28
$ws.Range("E77").Value = $cloudSyncFile
$ws.Range("C78").Value = $dateMar28
$ws.Range("E78").Value = $cloudSyncFile
$ws.Range("C79").Value = $dateMar28
$ws.Range("E79").Value = $cloudSyncFile

# Move the active selection, as a user would after finishing the review
$ws.Range("E70").Select() | Out-Null
